$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5 from 2023-09-06 (45175) to 2023-09-14 (45183)
$ws.Range("C2").Value = 45183
$ws.Range("C3").Value = 45183
$ws.Range("C4").Value = 45183
$ws.Range("C5").Value = 45183
